$wb = $excel.ActiveWorkbook

# --- Sheet 1: semantic_aspect_model_schema -------------------------------
# Widen column A (was squeezed to fit "id"; now needs room for "dtwin_id")
# and rename the digital-twin "id" header to "dtwin_id" so it no longer
# collides with the aspect model's own "id" column name.
$wsSchema = $wb.Worksheets.Item("semantic_aspect_model_schema")
$wsSchema.Columns.Item(1).ColumnWidth = 8.8
$wsSchema.Range("A1").Value = "dtwin_id"

# --- Sheet 2: description -------------------------------------------------
# Update the matching "Column Name" entry for the digital twin id field.
$wsDescription = $wb.Worksheets.Item("description")
$wsDescription.Range("A5").Value = "dtwin_id"

# --- Sheet 3: metadata (hidden) -------------------------------------------
# Bump the provenance info to the new upstream commit this template is
# regenerated from.
$wsMetadata = $wb.Worksheets.Item("metadata")
$wsMetadata.Range("B2").Value = "41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B3").Value = "https://github.com/dataspacesolutions/sldt-semantic-models/commit/41f43fae0e26ae5cfe94c2ce213309dcee6a0803"
$wsMetadata.Range("B4").Value = "2025-03-10 14:48:29+00:00"
